# Progetto linguaggi a, iniziata relazione
#
# The worksheet is an LR/predictive-parser trace table (columns: stack,
# input, action). Every row whose "action" column (C) was previously left
# blank actually represents an implicit "acc"(ept)/shift step; this edit
# fills those blanks in with the literal text "acc" (and, for the very
# last row of the trace - the terminating "$" / "$" row - with "halt"),
# formatted in a distinguishing green font so the newly-added annotations
# stand out from the original production-rule actions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose column-C cell was empty and now reads "acc".
$accRows = @(2, 4, 8, 11, 15, 17, 19, 22, 24, 26, 29, 31, 35, 39, 41, 44, 48, 51, 54, 57, 58)

# Final row's column-C cell was empty and now reads "halt".
$haltRow = 59

# Green font color used for the newly-inserted annotations (RGB 00B050,
# the same green Excel's "Font Color" gallery calls "Green"). COM/VBA
# colors are encoded little-endian BGR: R + G*256 + B*65536.
$accentGreen = 0 + (176 * 256) + (80 * 65536)

foreach ($r in $accRows) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = "acc"
    $cell.Font.Color = $accentGreen
}

$haltCell = $ws.Cells.Item($haltRow, 3)
$haltCell.Value = "halt"
$haltCell.Font.Color = $accentGreen

# Scroll the view down so row 49 is at the top (was row 28), matching
# where editing left off.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 49
$excel.ActiveWindow.ScrollColumn = 1

# Keep the active selection at A1, as in the original file.
$ws.Range("A1").Select()
